$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 184 ("「子供たちよ、いつでも思い出すんだよ」...") was removed from the
# source post list. Deleting the entire row shifts every following row
# (185-250) up by one, which matches the target diff (no content change,
# just renumbering), and updates the used range/dimension accordingly.
$ws.Rows.Item(184).Delete()
